$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Comment 0: zero-length range right before "Requirements Definition:"
# ---------------------------------------------------------------------------
$pReq = $d.Paragraphs(4).Range
$rC0 = $d.Range($pReq.Start, $pReq.Start)
$d.Comments.Add($rC0, "Thanks for mentioning Pythons  unittest framework to me in your comments. I took a quick look at it and it is perfect for the testing in this section.") | Out-Null
$d.Comments(1).Author = "Ray, Zach"
$d.Comments(1).Initial = "RZ"

# ---------------------------------------------------------------------------
# Comment 1: the word "database" inside the "A database" list item
# ---------------------------------------------------------------------------
$pDb = $d.Paragraphs(11).Range
$dbStart = $pDb.Start + 2
$dbEnd = $dbStart + 8
$rC1 = $d.Range($dbStart, $dbEnd)
$d.Comments.Add($rC1, "Are you able to store all info for the program in one database? Info for users, and info for courses?") | Out-Null
$d.Comments(2).Author = "Ray, Zach"
$d.Comments(2).Initial = "RZ"

# ---------------------------------------------------------------------------
# Comment 2: "System " inside "Integration and System Testing:"
# ---------------------------------------------------------------------------
$pInt = $d.Paragraphs(19).Range
$sysStart = $pInt.Start + 16
$sysEnd = $sysStart + 7
$rC2 = $d.Range($sysStart, $sysEnd)
$d.Comments.Add($rC2, "When do you plan on the full system being complete?") | Out-Null
$d.Comments(3).Author = "Ray, Zach"
$d.Comments(3).Initial = "RZ"
